# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计", populated
#   with the quarter's fund-holdings table (mirrors the "2021-Q4" sheet's
#   layout/columns, with a couple of extra columns).
# - Update the "总计" (totals) summary sheet with a new top row for 2022-Q1
#   (pushing the existing 2021-Q4 summary row down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4)
$newSheet.Name = "2022-Q1"

# Copy the header row's style (bold + bordered, centered) and the first
# column's style from the existing "2021-Q4" sheet so the new sheet matches
# the established look exactly.
$q4.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$q4.Range("A2").Copy()
$newSheet.Range("A2:A12").PasteSpecial(-4122)

# Header row
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Fund holding rows: code, name, fund size, total stock position,
# position ratio, holding market value (亿元), position rank
$fundRows = @(
    @("161724", "招商中证煤炭等权指数（LOF）", "21.53", "94.64", "3.25", "0.6997", 8),
    @("010779", "西部利得量化优选一年持有期混合A", "11.33", "87.57", "2.35", "0.2663", 4),
    @("501059", "西部利得中证国有企业红利指数增强（LOF）A", "2.22", "92.96", "2.88", "0.0639", 8),
    @("010780", "西部利得量化优选一年持有期混合C", "1.44", "87.57", "2.35", "0.0338", 4),
    @("009439", "西部利得中证国有企业红利指数增强（LOF）C", "0.98", "92.96", "2.88", "0.0282", 8),
    @("004352", "北信瑞丰研究精选股票", "0.49", "92.71", "1.13", "0.0055", 5),
    @("009263", "华宝红利精选混合A", "0.46", "83.67", "0.98", "0.0045", 7),
    @("007808", "北信瑞丰量化优选灵活配置混合", "0.24", "79.84", "1.09", "0.0026", 3),
    @("006857", "蜂巢卓睿灵活配置混合A", "0.14", "78.16", "1.17", "0.0016", 10),
    @("010841", "华宝红利精选混合C", "0.16", "83.67", "0.98", "0.0016", 7),
    @("006858", "蜂巢卓睿灵活配置混合C", "0.04", "78.16", "1.17", "0.0005", 10)
)

$r = 2
foreach ($fund in $fundRows) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2
    # Fund code and the percentage-style columns are stored as plain text in
    # the source data (e.g. fund codes keep their leading zeros), so prefix
    # with an apostrophe to force text just like typing them in by hand.
    $newSheet.Cells.Item($r, 2).Value = "'" + $fund[0]
    $newSheet.Cells.Item($r, 3).Value = $fund[1]
    $newSheet.Cells.Item($r, 4).Value = "'" + $fund[2]
    $newSheet.Cells.Item($r, 5).Value = "'" + $fund[3]
    $newSheet.Cells.Item($r, 6).Value = "'" + $fund[4]
    $newSheet.Cells.Item($r, 7).Value = "'" + $fund[5]
    $newSheet.Cells.Item($r, 8).Value = $fund[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: push the existing 2021-Q4 row down to row 3 and
#    insert the new 2022-Q1 totals at row 2.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Preserve row 3's "index" cell style (A column) to match row 2's, then move
# the existing 2021-Q4 totals down one row.
$totalSheet.Cells.Item(2,1).Copy()
$totalSheet.Cells.Item(3,1).PasteSpecial(-4122)

$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(3,2).Value = "2021-Q4"
$totalSheet.Cells.Item(3,3).Value = 1
$totalSheet.Cells.Item(3,4).Value = 0.9

$totalSheet.Cells.Item(2,2).Value = "2022-Q1"
$totalSheet.Cells.Item(2,3).Value = 11
$totalSheet.Cells.Item(2,4).Value = 1.11

Write-Output "2022-Q1 sheet added and 总计 sheet updated"
